# overall offline analysis, maxwellian init option
#
# Summary of the edits reproduced below:
#   1) Rename sheet "Arcs" -> "Arc0" (same sheetId / position).
#   2) Cameras!C6: camera opening angle 88.3 -> 88.0172525718237.
#   3) Cameras sheet: move the (bottom-pane) selection from C24 to B15.
#   4) Sim sheet: scroll the frozen window up a little (was parked at
#      row 35, now parked at row 32); it also stops being the front tab.
#   5) Arc0 becomes the new front/active tab, and its window scrolls back
#      to the left (was parked at column W, now back at column A).

$wb = $excel.ActiveWorkbook

$simSheet     = $wb.Worksheets.Item("Sim")
$camerasSheet = $wb.Worksheets.Item("Cameras")
$arcsSheet    = $wb.Worksheets.Item("Arcs")

# --- Sim: nudge the frozen-pane scroll position up to row 32 -------------
$simSheet.Activate()
$excel.ActiveWindow.ScrollRow = 32

# --- Cameras: update the camera opening angle and move the selection -----
$camerasSheet.Activate()
$camerasSheet.Range("C6").Value = 88.0172525718237
$camerasSheet.Range("B15").Select()

# --- Rename "Arcs" -> "Arc0" ----------------------------------------------
$arcsSheet.Name = "Arc0"

# --- Arc0 becomes the active sheet, scrolled back to column A ------------
$arcsSheet.Activate()
$excel.ActiveWindow.ScrollColumn = 1
